$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '89.252.85'
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('D3').Value = '3.158.53'
$ws.Range('E3').Value = '  -3.71%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.42'
$ws.Range('E5').Value = '  -0.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '612.68'
$ws.Range('E6').Value = '  -2.56%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.385'
$ws.Range('E7').Value = '  +2.40%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.686'
$ws.Range('E8').Value = '  -5.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').Value = '3.161.83'
$ws.Range('E10').Value = '  -3.37%  '
$ws.Range('E11').Value = '  -0.83%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.176'
$ws.Range('E12').Value = '  -6.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000251'
$ws.Range('E13').Value = '  -4.75%  '
$ws.Range('D14').Value = '89.316.32'
$ws.Range('E14').Value = '  +1.12%  '
$ws.Range('D15').Value = '3.748.33'
$ws.Range('E15').Value = '  -3.46%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '32.56'
$ws.Range('E16').Value = '  -4.78%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.21'
$ws.Range('E17').Value = '  -5.23%  '
$ws.Range('D18').Value = '3.153.71'
$ws.Range('E18').Value = '  -4.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.27'
$ws.Range('E19').Value = '  +3.16%  '
$ws.Range('E20').Value = '  -5.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '434.20'
$ws.Range('E21').Value = '  -0.67%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.55'
$ws.Range('E22').Value = '  -4.01%  '
$ws.Range('B23').Value = 'PEPE'
$ws.Range('C23').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0000186'
$ws.Range('E23').Value = '  +36.88%  '
$ws.Range('E24').Value = '  -5.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.09'
$ws.Range('E25').Value = '  -3.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.64'
$ws.Range('E26').Value = '  -5.30%  '
$ws.Range('D27').Value = '3.342.10'
$ws.Range('E27').Value = '  -3.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '74.95'
$ws.Range('E28').Value = '  -2.83%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.166'
$ws.Range('E30').Value = '  -7.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.08'
$ws.Range('E32').Value = '  +31.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '8.38'
$ws.Range('E33').Value = '  -4.47%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '527.90'
$ws.Range('E34').Value = '  -7.73%  '
$ws.Range('E35').Value = '  -2.20%  '
$ws.Range('E36').Value = '  -6.32%  '
$ws.Range('E37').Value = '  -8.89%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '22.30'
$ws.Range('E38').Value = '  +2.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '21.74'
$ws.Range('E39').Value = '  -4.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.19%  '
$ws.Range('E41').Value = '  -9.85%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.91'
$ws.Range('E43').Value = '  -5.92%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.370'
$ws.Range('E44').Value = '  -7.95%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '149.54'
$ws.Range('E45').Value = '  -2.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '43.88'
$ws.Range('E46').Value = '  -2.12%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '171.13'
$ws.Range('E47').Value = '  -5.42%  '
$ws.Range('E48').Value = '  -10.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.23'
$ws.Range('E49').Value = '  -7.29%  '
$ws.Range('B50').Value = 'Filecoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.03'
$ws.Range('E50').Value = '  -4.59%  '
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.605'
$ws.Range('E51').Value = '  -3.60%  '
